# Applies the scheduled-runner cell value updates to the Sheets workbook.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds crafting-profit
# data; columns H-N (currentAveragePrice.. LeveProfitHQ) are refreshed per row.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1036.6
$ws.Range("I6").Value = 444.91666
$ws.Range("J6").Value = 3403.3333
$ws.Range("K6").Value = 1334.74998
$ws.Range("L6").Value = 10209.9999
$ws.Range("M6").Value = -1222.74998
$ws.Range("N6").Value = -10433.9999
$ws.Range("H7").Value = 250
$ws.Range("I7").Value = 250
$ws.Range("K7").Value = 250
$ws.Range("M7").Value = -138
$ws.Range("H14").Value = 250
$ws.Range("I14").Value = 250
$ws.Range("K14").Value = 250
$ws.Range("M14").Value = -59
$ws.Range("H31").Value = 251.83333
$ws.Range("I31").Value = 102.4
$ws.Range("J31").Value = 999
$ws.Range("K31").Value = 307.2
$ws.Range("L31").Value = 2997
$ws.Range("M31").Value = -77.20000000000005
$ws.Range("N31").Value = -3457
$ws.Range("H38").Value = 11034.667
$ws.Range("J38").Value = 26193.285
$ws.Range("L38").Value = 78579.855
$ws.Range("N38").Value = -79323.855
$ws.Range("H64").Value = 12861571
$ws.Range("I64").Value = 45002500
$ws.Range("J64").Value = 5199.8
$ws.Range("K64").Value = 45002500
$ws.Range("L64").Value = 5199.8
$ws.Range("M64").Value = -45002252
$ws.Range("N64").Value = -5695.8
$ws.Range("H67").Value = 12861571
$ws.Range("I67").Value = 45002500
$ws.Range("J67").Value = 5199.8
$ws.Range("K67").Value = 45002500
$ws.Range("L67").Value = 5199.8
$ws.Range("M67").Value = -45001642
$ws.Range("N67").Value = -6915.8
$ws.Range("H125").Value = 2345.8
$ws.Range("I125").Value = 2133
$ws.Range("K125").Value = 19197
$ws.Range("M125").Value = -16737
$ws.Range("H135").Value = 1871.25
$ws.Range("I135").Value = 1871.25
$ws.Range("K135").Value = 16841.25
$ws.Range("M135").Value = -14306.25
$ws.Range("H137").Value = 1733.5
$ws.Range("I137").Value = 1733.5
$ws.Range("K137").Value = 5200.5
$ws.Range("M137").Value = -2650.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents() | Out-Null
$ws.Range("H32").Value = 998.5
$ws.Range("I32").Value = 747.5789
$ws.Range("K32").Value = 747.5789
$ws.Range("M32").Value = -460.5789
$ws.Range("H38").Value = 1000
$ws.Range("I38").Value = 1000
$ws.Range("K38").Value = 1000
$ws.Range("M38").Value = -533
$ws.Range("H42").Value = 20031
$ws.Range("J42").Value = 20031
$ws.Range("L42").Value = 20031
$ws.Range("N42").Value = -21003
$ws.Range("H97").Value = 551.8
$ws.Range("I97").Value = 551.8
$ws.Range("K97").Value = 551.8
$ws.Range("M97").Value = -55.79999999999995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents() | Out-Null
$ws.Range("H22").Value = 212.33333
$ws.Range("I22").Value = 212.33333
$ws.Range("K22").Value = 212.33333
$ws.Range("M22").Value = -39.33332999999999
$ws.Range("H95").Value = 11222
$ws.Range("J95").Value = 12810.75
$ws.Range("L95").Value = 12810.75
$ws.Range("N95").Value = -18302.75
$ws.Range("H105").Value = 5498.2856
$ws.Range("I105").Value = 5664.6665
$ws.Range("K105").Value = 5664.6665
$ws.Range("M105").Value = -3917.6665
$ws.Range("H135").Value = 37500
$ws.Range("J135").Value = 37500
$ws.Range("L135").Value = 37500
$ws.Range("N135").Value = -47640

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1617
$ws.Range("I31").Value = 1500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1205
$ws.Range("H34").Value = 1617
$ws.Range("I34").Value = 1500
$ws.Range("K34").Value = 1500
$ws.Range("M34").Value = -1298
$ws.Range("H62").Value = 56259496
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents() | Out-Null
$ws.Range("H65").Value = 56259496
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents() | Out-Null
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents() | Out-Null

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1005.2
$ws.Range("I7").Value = 1253.75
$ws.Range("K7").Value = 3761.25
$ws.Range("M7").Value = -3649.25
$ws.Range("H75").Value = 600
$ws.Range("J75").Value = 600
$ws.Range("L75").Value = 1800
$ws.Range("N75").Value = -3796
$ws.Range("H78").Value = 600
$ws.Range("J78").Value = 600
$ws.Range("L78").Value = 5400
$ws.Range("N78").Value = -15384
$ws.Range("H113").Value = 591.61536
$ws.Range("J113").Value = 649.3
$ws.Range("L113").Value = 1947.9
$ws.Range("N113").Value = -6287.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 417416.5
$ws.Range("I3").Value = 834000
$ws.Range("J3").Value = 833
$ws.Range("K3").Value = 834000
$ws.Range("L3").Value = 833
$ws.Range("M3").Value = -833884
$ws.Range("N3").Value = -1065
$ws.Range("H70").Value = 99998.5
$ws.Range("J70").Value = 99998
$ws.Range("L70").Value = 99998
$ws.Range("N70").Value = -100538
$ws.Range("H73").Value = 99998.5
$ws.Range("J73").Value = 99998
$ws.Range("L73").Value = 99998
$ws.Range("N73").Value = -101870
$ws.Range("H102").Value = 8331.666999999999
$ws.Range("I102").Value = 8331.666999999999
$ws.Range("K102").Value = 8331.666999999999
$ws.Range("M102").Value = -6709.666999999999
$ws.Range("H132").Value = 5891.8184
$ws.Range("I132").Value = 5312.222
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 15936.666
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -13406.666
$ws.Range("N132").Value = -30560

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4527.273
$ws.Range("I7").Value = 4480
$ws.Range("K7").Value = 4480
$ws.Range("M7").Value = -4368
$ws.Range("H33").Value = 352533.34
$ws.Range("I33").Value = 1000000
$ws.Range("J33").Value = 28800
$ws.Range("K33").Value = 1000000
$ws.Range("L33").Value = 28800
$ws.Range("M33").Value = -999710
$ws.Range("N33").Value = -29380
$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 50000
$ws.Range("K74").Value = 50000
$ws.Range("M74").Value = -49002
$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 50000
$ws.Range("K77").Value = 150000
$ws.Range("M77").Value = -145008
$ws.Range("H93").Value = 15066.333
$ws.Range("I93").Value = 15066.333
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 15066.333
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -13818.333
$ws.Range("N93").ClearContents() | Out-Null
$ws.Range("H126").Value = 4527.273
$ws.Range("I126").Value = 4480
$ws.Range("K126").Value = 13440
$ws.Range("M126").Value = -10970
$ws.Range("H132").Value = 4589
$ws.Range("J132").Value = 4247.5
$ws.Range("L132").Value = 12742.5
$ws.Range("N132").Value = -17802.5
$ws.Range("H136").Value = 6174
$ws.Range("I136").Value = 5398.6665
$ws.Range("K136").Value = 16195.9995
$ws.Range("M136").Value = -13645.9995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 2239.2
$ws.Range("I6").Value = 196.66667
$ws.Range("J6").Value = 5303
$ws.Range("K6").Value = 196.66667
$ws.Range("L6").Value = 5303
$ws.Range("M6").Value = -81.66667000000001
$ws.Range("N6").Value = -5533
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents() | Out-Null
$ws.Range("N62").ClearContents() | Out-Null
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents() | Out-Null
$ws.Range("N65").ClearContents() | Out-Null
$ws.Range("H122").Value = 4030.75
$ws.Range("I122").Value = 4030.75
$ws.Range("K122").Value = 12092.25
$ws.Range("M122").Value = -9642.25
$ws.Range("H131").Value = 69994
$ws.Range("J131").Value = 69994
$ws.Range("L131").Value = 69994
$ws.Range("N131").Value = -80074
$ws.Range("H136").Value = 6968.3335
$ws.Range("I136").Value = 7000
$ws.Range("K136").Value = 21000
$ws.Range("M136").Value = -18450

